$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This script re-applies yellow highlighting to the "(15 points)" Sweepstakes
# user-story bullet block (the paragraph itself plus its four sub-bullets),
# merges the two runs that used to read "(15 points" + ")" into a single
# "(15 points)" run, merges the two trailing runs of that same paragraph into
# one highlighted run, and moves the stray "_GoBack" bookmark from the end of
# the "RegisterContestant" bullet down to the end of the "PickWinner" bullet.
#
# We rebuild each affected paragraph's XML in place via Range.InsertXML so
# that the paragraph-mark run properties (w:pPr/w:rPr) pick up the highlight
# too (simple Range.HighlightColorIndex assignments only ever touch the run
# content, never the paragraph mark's own formatting).
# ---------------------------------------------------------------------------

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Paragraph: "(15 points) As a developer, I want to create a Sweepstakes
# class that uses the Dictionary data structure as an underlying structure.
# The Sweepstakes class will have the following methods with full
# implementation (write the functionality) of each method:"
$p7 = $d.Paragraphs.Item(7)
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="54AA6FD7" w14:textId="047C9DB3" w:rsidR="3A99A137" w:rsidRPr="00ED6A79" w:rsidRDefault="3A99A137" w:rsidP="00ED6A79">
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>(15 points)</w:t>
  </w:r>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve"> As a developer, I want to create a Sweepstakes class that uses the Dictionary data structure as an underlying structure. The Sweepstakes class will have the following methods with full implementation (write the functionality) of each method:</w:t>
  </w:r>
</w:p>
'@
$p7.Range.InsertXML($xml7)

# Paragraph: "Sweepstakes(string name)"
$p8 = $d.Paragraphs.Item(8)
$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="08A89551" w14:textId="77777777" w:rsidR="002D4749" w:rsidRPr="00ED6A79" w:rsidRDefault="002D4749" w:rsidP="00ED6A79">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>Sweepstakes(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>string name)</w:t>
  </w:r>
</w:p>
'@
$p8.Range.InsertXML($xml8)

# Paragraph: "void RegisterContestant(Contestant contestant)" - loses the
# "_GoBack" bookmark (it moves to the next bullet below).
$p9 = $d.Paragraphs.Item(9)
$xml9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="46750CE2" w14:textId="45ADEED3" w:rsidR="002D4749" w:rsidRPr="00ED6A79" w:rsidRDefault="002D4749" w:rsidP="00ED6A79">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve">void </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>RegisterContestant</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>Contestant contestant)</w:t>
  </w:r>
</w:p>
'@
$p9.Range.InsertXML($xml9)

# Paragraph: "string PickWinner()" - gains the "_GoBack" bookmark.
$p10 = $d.Paragraphs.Item(10)
$xml10 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="137C0662" w14:textId="3CAFEC18" w:rsidR="002D4749" w:rsidRPr="00ED6A79" w:rsidRDefault="002D4749" w:rsidP="00ED6A79">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve">string </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>PickWinner</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@
$p10.Range.InsertXML($xml10)

# Paragraph: "void PrintContestantInfo(Contestant contestant)"
$p11 = $d.Paragraphs.Item(11)
$xml11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6E872AAF" w14:textId="4813AB11" w:rsidR="002D4749" w:rsidRDefault="002D4749" w:rsidP="00ED6A79">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t xml:space="preserve">void </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>PrintContestantInfo</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r w:rsidRPr="00ED6A79">
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>Contestant contestant)</w:t>
  </w:r>
</w:p>
'@
$p11.Range.InsertXML($xml11)

Write-Host "Sweepstakes class bullets re-highlighted."
